$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Column L corresponds to site "US10" (Cincinnati). Mark the following
# treatment components as present (1) for that site.
$rows = @(2, 4, 12, 15, 17, 18, 19, 20, 21, 22, 23, 25, 26, 27, 32, 37, 38)

foreach ($r in $rows) {
    $ws.Range("L$r").Value = 1
}

# Update the active selection to match the edited workbook state.
$ws.Range("L39").Select()
